# Update NATMI Gdf7-Acvr2a LR-pair sheet with refreshed TPM-derived values.
# Rows 2-5 (Sending cluster = MuSCs) get updated edge-weight metrics, and
# four new rows 6-9 (Sending cluster = Resolving-Mac) are appended for the
# same four target clusters (ECs, FAPs, MuSCs, Resolving-Mac).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Gdf7"
$ws.Range("C2").Value = "Acvr2a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02913733333333333
$ws.Range("H2").Value = 0.087412
$ws.Range("I2").Value = 0.2991461473965196
$ws.Range("J2").Value = 0.2991461473965196
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 0.4746019194337778
$ws.Range("R2").Value = 4.271417274904
$ws.Range("S2").Value = 0.06512126682555208
$ws.Range("T2").Value = 0.06512126682555208

# Row 3
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Gdf7"
$ws.Range("C3").Value = "Acvr2a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02913733333333333
$ws.Range("H3").Value = 0.087412
$ws.Range("I3").Value = 0.2991461473965196
$ws.Range("J3").Value = 0.2991461473965196
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83272
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 0.8045081911822222
$ws.Range("R3").Value = 7.24057372064
$ws.Range("S3").Value = 0.1103884970457435
$ws.Range("T3").Value = 0.1103884970457435

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Gdf7"
$ws.Range("C4").Value = "Acvr2a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02913733333333333
$ws.Range("H4").Value = 0.087412
$ws.Range("I4").Value = 0.2991461473965196
$ws.Range("J4").Value = 0.2991461473965196
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 0.7653306669666666
$ws.Range("R4").Value = 6.8879760027
$ws.Range("S4").Value = 0.1050128550528718
$ws.Range("T4").Value = 0.1050128550528718

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gdf7"
$ws.Range("C5").Value = "Acvr2a"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02913733333333333
$ws.Range("H5").Value = 0.087412
$ws.Range("I5").Value = 0.2991461473965196
$ws.Range("J5").Value = 0.2991461473965196
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 0.1357277398071111
$ws.Range("R5").Value = 1.221549658264
$ws.Range("S5").Value = 0.01862352847235225
$ws.Range("T5").Value = 0.01862352847235225

# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Gdf7"
$ws.Range("C6").Value = "Acvr2a"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.06826433333333333
$ws.Range("H6").Value = 0.204793
$ws.Range("I6").Value = 0.7008538526034804
$ws.Range("J6").Value = 0.7008538526034804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.28844733333333
$ws.Range("N6").Value = 48.865342
$ws.Range("O6").Value = 0.2176904746803693
$ws.Range("P6").Value = 0.2176904746803693
$ws.Range("Q6").Value = 1.111919998245111
$ws.Range("R6").Value = 10.007279984206
$ws.Range("S6").Value = 0.1525692078548173
$ws.Range("T6").Value = 0.1525692078548173

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Gdf7"
$ws.Range("C7").Value = "Acvr2a"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.06826433333333333
$ws.Range("H7").Value = 0.204793
$ws.Range("I7").Value = 0.7008538526034804
$ws.Range("J7").Value = 0.7008538526034804
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.61090666666666
$ws.Range("N7").Value = 82.83272
$ws.Range("O7").Value = 0.3690119294748028
$ws.Range("P7").Value = 0.3690119294748029
$ws.Range("Q7").Value = 1.884840136328889
$ws.Range("R7").Value = 16.96356122696
$ws.Range("S7").Value = 0.2586234324290594
$ws.Range("T7").Value = 0.2586234324290594

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Gdf7"
$ws.Range("C8").Value = "Acvr2a"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.06826433333333333
$ws.Range("H8").Value = 0.204793
$ws.Range("I8").Value = 0.7008538526034804
$ws.Range("J8").Value = 0.7008538526034804
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.266325
$ws.Range("N8").Value = 78.798975
$ws.Range("O8").Value = 0.3510419771967738
$ws.Range("P8").Value = 0.3510419771967739
$ws.Range("Q8").Value = 1.793053165241666
$ws.Range("R8").Value = 16.137478487175
$ws.Range("S8").Value = 0.2460291221439021
$ws.Range("T8").Value = 0.2460291221439021

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Gdf7"
$ws.Range("C9").Value = "Acvr2a"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.06826433333333333
$ws.Range("H9").Value = 0.204793
$ws.Range("I9").Value = 0.7008538526034804
$ws.Range("J9").Value = 0.7008538526034804
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.658207333333333
$ws.Range("N9").Value = 13.974622
$ws.Range("O9").Value = 0.06225561864805391
$ws.Range("P9").Value = 0.06225561864805392
$ws.Range("Q9").Value = 0.3179894181384444
$ws.Range("R9").Value = 2.861904763246
$ws.Range("S9").Value = 0.04363209017570167
$ws.Range("T9").Value = 0.04363209017570167
